# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.935.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.765.52"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "402.82"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.40"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.754.14"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -6.18%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.166"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -10.18%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -12.18%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.367.65"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.67"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.57"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +12.48%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.748.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.37"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.185.15"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "411.55"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -8.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.40"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -8.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.70"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.03"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.72"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +14.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "36.12"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.09"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -7.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.31"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -8.84%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.73"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.31"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.153"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "38.82"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -7.84%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.94"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0728"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.30%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.88"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.135"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -8.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.90"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.16"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +19.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "144.97"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.24"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.07%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.22"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.291"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.53%  "
